# Generate Report for Handoff
# Update the "Latest Handoff" timestamps for the 9ec7f35c-8c26-4081-a14f-1ff8ed6414dc.md
# entry across the Overview, zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-03-24 06:44:43"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-24 06:44:39"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-24 06:44:43"
